# Refresh the crypto price table: Price (column D) and Volume(1h) (column E)
# values are updated to the latest pulled figures, per the automated
# "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as "1.729.42" or "0.9999" that must stay
# literal text (thousand-dot formatting, fixed trailing zeros, etc.). Mark the
# column as Text before writing so Excel does not reinterpret the values as
# numbers, then restore the default "Normal" style once done so no visible
# formatting changes are left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.548.02"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.729.12"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "245.21"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "0.4801"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "0.2666"
$ws.Range("E8").Value = "  -1.62%  "
$ws.Range("D9").Value = "0.06178"
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").Value = "1.729.05"
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("D11").Value = "0.07179"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").Value = "15.55"
$ws.Range("E12").Value = "  -2.66%  "
$ws.Range("D13").Value = "0.6086"
$ws.Range("E13").Value = "  -2.42%  "
$ws.Range("D14").Value = "4.529"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "77.21"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "26.552.93"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "0.000006959"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").Value = "1.953.86"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("E22").Value = "  -2.81%  "
$ws.Range("D23").Value = "8.785"
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").Value = "5.229"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").Value = "137.05"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("E27").Value = "  -3.56%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "107.49"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("D32").Value = "3.683"
$ws.Range("E32").Value = "  -2.37%  "
$ws.Range("D33").Value = "0.04501"
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("D34").Value = "2.616"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("D35").Value = "1.001"
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("D36").Value = "0.6308"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").Value = "0.9076"
$ws.Range("E37").Value = "  -4.67%  "
$ws.Range("D38").Value = "2.046"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("D39").Value = "2.400"
$ws.Range("E39").Value = "  -4.44%  "
$ws.Range("D40").Value = "1.004"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D41").Value = "103.41"
$ws.Range("E41").Value = "  -9.84%  "
$ws.Range("D42").Value = "0.01502"
$ws.Range("E42").Value = "  -0.74%  "
$ws.Range("D43").Value = "5.493"
$ws.Range("E43").Value = "  -4.28%  "
$ws.Range("D44").Value = "0.3888"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("E45").Value = "  +3.80%  "
$ws.Range("E46").Value = "  -2.15%  "
$ws.Range("D47").Value = "0.05384"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("D48").Value = "30.67"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("D49").Value = "7.829"
$ws.Range("E49").Value = "  -2.48%  "
$ws.Range("D50").Value = "1.247"
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").Value = "0.3410"
$ws.Range("E51").Value = "  -1.68%  "

# Restore the original (default) cell style on column D now that the text
# values are safely in place.
$ws.Range("D2:D51").Style = "Normal"
